# Update Moving File To Outside Directory
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 8 data updates
$ws.Range("A8").Value = "RUN"
$ws.Range("M8").Value = "CLM2023010097"
$ws.Range("N8").Value = "04/08/2022"
$ws.Range("O8").Value = "16/09/2022"

# Update the view state (scroll position / active selection)
$ws.Application.Goto($ws.Range("A1"), $true)
$ws.Range("J5").Select()
$excel.ActiveWindow.ScrollRow = 5
$excel.ActiveWindow.ScrollColumn = 10
$ws.Range("P8").Select()
